$wb = $excel.ActiveWorkbook

# --- "shrub locations" sheet: update selection, no longer the active tab ---
$wsLocations = $wb.Worksheets.Item("shrub locations")
$wsLocations.Activate()
$wsLocations.Range("D2:E31").Select()

# --- "shrubs" sheet: update selection + page setup ---
$wsShrubs = $wb.Worksheets.Item("shrubs")
$wsShrubs.Activate()
$wsShrubs.Range("D3:F32").Select()
$wsShrubs.PageSetup.Orientation = 1

# --- "annuals" sheet: remove the merged "species and abundance" header row,
#     shifting everything up by one, then set it as the active tab with the
#     appropriate view/selection (whole used columns A:N selected) ---
$wsAnnuals = $wb.Worksheets.Item("annuals")
$wsAnnuals.Rows.Item(1).Delete()
$wsAnnuals.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$wsAnnuals.Range("A1:N1048576").Select()
